$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New station codes (Column A), descriptions (Column B), and values (Column C)
$codes = @("HarveyCanalNorth", "HarveyCanalBoom", "BayouBienv", "BaraPass", "FreshCanal", "CalcRiv")
$descs = @(
    "Harvey Canal Sector Gates North / Prot Side nr Lapalco Blvd",
    "Harvey Canal at Boomtown Casion",
    "Bayou Bienvenue Floodgate",
    "Barataria Pass at Grand Isle",
    "Freshwater Canal at Freshwater Bayou Lock South",
    "Calcasieu River at Cameron"
)
$values = @(2, 0, 3, 0, 0, 0)

$startRow = 17
$fmtA = $ws.Cells.Item(1, 1).NumberFormat
$fmtC = $ws.Cells.Item(1, 3).NumberFormat

# Fill column A first (matches shared-string insertion order of the source workbook)
for ($i = 0; $i -lt $codes.Count; $i++) {
    $r = $startRow + $i
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $codes[$i]
    $cell.NumberFormat = $fmtA
}

# Then fill column B
for ($i = 0; $i -lt $descs.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value2 = $descs[$i]
}

# Then fill column C
for ($i = 0; $i -lt $values.Count; $i++) {
    $r = $startRow + $i
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $values[$i]
    $cell.NumberFormat = $fmtC
}

$ws.Range("C23").Select()

$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

$wb.Save()
